$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# --- Header row (row 1) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# --- Data row (row 2) ---
$ws.Range("A2").Value = 71
$ws.Range("B2").Value = "富達全聚焦"
$ws.Range("C2").Value = "楊際英"
$ws.Range("D2").Value = "台北富邦商業銀行金華分行"
$ws.Range("E2").Value = 1089.68
$ws.Range("F2").Value = 42715.46
$ws.Range("G2").Value = "美金"
$ws.Range("H2").Value = 1257756
$ws.Range("I2").Value = "fund"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "2012-04-20"
$ws.Range("L2").Value = "蔣乃辛"
$ws.Range("M2").Value = 1722
$ws.Range("N2").Value = "tmp7091"
$ws.Range("O2").Value = 71
